$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new mapped field row: id_platoon (new table field with no "Antigo" counterpart)
$ws.Range("B21").Value = "id_platoon"

# Fix typo: "banck_agency_conta" -> "bank_agency_conta"
$ws.Range("B11").Value = "bank_agency_conta"

# Leave selection on the cell that was last edited
$ws.Range("B11").Select() | Out-Null
